$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G:K on the new rows hold numeric-looking text ("0", "1", "75.00", "-", ...)
# that must stay text (matching the source sheet, which stores every cell as a
# string). Pre-format as Text so Excel's Value setter doesn't coerce them to numbers.
$ws.Range("G5:K7").NumberFormat = "@"

# Row 5 - duplicate of row 2's match (Sunrisers Hyderabad vs Chennai Super Kings)
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 13 2020"
$ws.Range("C5").Value = "Super Kings won by 20 runs"
$ws.Range("D5").Value = "Sunrisers Hyderabad"
$ws.Range("E5").Value = "Chennai Super Kings"
$ws.Range("F5").Value = "T Natarajan "
$ws.Range("G5").Value = "0"
$ws.Range("H5").Value = "1"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "0.00"

# Row 6 - duplicate of row 4's match (Sunrisers Hyderabad vs Royal Challengers Bangalore)
$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " September 21 2020"
$ws.Range("C6").Value = "RCB won by 10 runs"
$ws.Range("D6").Value = "Sunrisers Hyderabad"
$ws.Range("E6").Value = "Royal Challengers Bangalore"
$ws.Range("F6").Value = "T Natarajan "
$ws.Range("G6").Value = "3"
$ws.Range("H6").Value = "4"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "75.00"

# Row 7 - duplicate of row 3's match (Sunrisers Hyderabad vs Kings XI Punjab)
$ws.Range("A7").Value = " Dubai (DSC)"
$ws.Range("B7").Value = " October 24 2020"
$ws.Range("C7").Value = "Kings XI won by 12 runs"
$ws.Range("D7").Value = "Sunrisers Hyderabad"
$ws.Range("E7").Value = "Kings XI Punjab"
$ws.Range("F7").Value = "T Natarajan "
$ws.Range("G7").Value = "0"
$ws.Range("H7").Value = "0"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "-"
